$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value that "looks like" a number/currency
# (e.g. "$1,409.00") without letting Excel auto-convert it to a numeric
# cell. A leading apostrophe forces text entry (quote-prefixed); resetting
# the style back to Normal afterwards drops the quote-prefix formatting so
# the cell ends up as a plain shared-string cell, matching the existing
# "$1,409.00" cells elsewhere on this sheet (e.g. C18/C19).
function Set-TextValue($rng, $text) {
    $rng.Formula = "'" + $text
    $rng.Style = "Normal"
}

# New "Rewrite Transaction" style Issuance/Change comparison block added
# to the right of the existing data (mirrors the existing batch-process
# rows lower on the sheet).
$ws.Range("L10").Value = "Issuance"
Set-TextValue $ws.Range("M10") "$1,409.00"
$ws.Range("N10").Value = "Issuance"

$ws.Range("E12").Value = "Issuance"
Set-TextValue $ws.Range("F12") "$1,409.00"

$ws.Range("L12").Value = "Change"
Set-TextValue $ws.Range("M12") "$1,409.00"
$ws.Range("N12").Value = "Change"

# Leave the selection where the author left it after editing.
$ws.Range("S7").Select() | Out-Null
